$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.116.16'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.828.03'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.07'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4579'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +7.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3741'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07320'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8636'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.98'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '1.826.33'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.723'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.372'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.31%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.99'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +5.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07098'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008863'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.04'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '27.098.65'
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.205'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.19'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.232'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.53'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.290'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08913'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.201'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7654'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("E33").Value = '  +5.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.486'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01976'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05300'
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5381'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +7.01%  '
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.878'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("E42").Value = '  +2.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5225'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +11.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.640'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.74'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.999'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +11.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.25'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.690'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9259'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.41%  '
